# Apply "Add data for 2021-10-28" update to the carjacking-by-neighborhood-by-month workbook.
# This shifts the reporting window from "through October 19" to "through October 20"
# and updates the affected cell values accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab to reflect the new reporting date.
$ws.Name = "Through 2021-10-20"

# Update the column B header text (shared string) to match the new window label.
$ws.Range("B1").Value = "October 2021 (through October 20)"

# Update the individual data cells that changed.
$ws.Range("L2").Value = 13    # Garfield Park / October 2020
$ws.Range("B3").Value = 11    # Austin / October 2021 (through October 20)
$ws.Range("V3").Value = 4     # Austin / October 2019
$ws.Range("AP3").Value = 6    # Austin / October 2017
$ws.Range("L6").Value = 1     # Auburn Gresham / October 2020
$ws.Range("BJ8").Value = 1    # Little Italy, UIC / October 2015
$ws.Range("AF9").Value = 3    # Humboldt Park / October 2018
$ws.Range("L13").Value = 3    # South Shore / October 2020
$ws.Range("B16").Value = 1    # Avalon Park / October 2021 (through October 20)
$ws.Range("V18").Value = 1    # Little Village / October 2019
$ws.Range("E27").Value = 4    # West Loop / July 2021
$ws.Range("B29").Value = 2    # Archer Heights / October 2021 (through October 20)
$ws.Range("AZ33").Value = 1   # Hyde Park / October 2016
$ws.Range("AP38").Value = 4   # Englewood / October 2017
$ws.Range("AZ38").Value = 4   # Englewood / October 2016
$ws.Range("V41").Value = 2    # Washington Heights / October 2019
